$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.705.72"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "3.820.28"
$ws.Range("E3").Value = "  +2.17%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.75%  "

$ws.Range("D7").Value = "3.822.28"
$ws.Range("E7").Value = "  +2.29%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.43"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.477"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.83"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000252"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.10%  "

$ws.Range("D15").Value = "4.464.31"
$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").Value = "3.821.86"
$ws.Range("E16").Value = "  +2.16%  "

$ws.Range("D17").Value = "69.803.41"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.16%  "

$ws.Range("E19").Value = "  -3.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "503.68"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.45"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.732"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.46"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000142"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.40"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.74%  "

$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.97"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.79"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.66%  "

$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("E36").Value = "  -2.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.04"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.83%  "

$ws.Range("E38").Value = "  +2.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "487.31"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +14.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.334"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.04"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.67%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.70"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.46%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.97"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.20"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.48"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.16%  "

$ws.Range("D46").Value = "2.916.04"
$ws.Range("E46").Value = "  -3.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0359"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "139.88"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.21%  "

$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.75"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.41"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.60%  "
